$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column I values
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 8

# New column J values
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 4
$ws.Range("J4").Value = 4
$ws.Range("J5").Value = 2
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 8
